$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1584
$ws.Range("F3").Value = 121
$ws.Range("F5").Value = 9043
$ws.Range("F6").Value = 258
$ws.Range("F7").Value = 113
$ws.Range("F8").Value = 1262
$ws.Range("F10").Value = 586
$ws.Range("F11").Value = 598
$ws.Range("F13").Value = 134
$ws.Range("F14").Value = 289
$ws.Range("F16").Value = 53
$ws.Range("F17").Value = 1481
$ws.Range("F18").Value = 1317
$ws.Range("F21").Value = 1364
$ws.Range("F22").Value = 77
$ws.Range("F23").Value = 225
$ws.Range("F25").Value = 86
$ws.Range("F26").Value = 51
$ws.Range("F27").Value = 63
$ws.Range("F28").Value = 298
$ws.Range("F29").Value = 298
$ws.Range("F30").Value = 1065
$ws.Range("F31").Value = 10
$ws.Range("F33").Value = 225
$ws.Range("F34").Value = 191
$ws.Range("F42").Value = 21
$ws.Range("F43").Value = 487
$ws.Range("F44").Value = 1228
$ws.Range("F45").Value = 683
$ws.Range("F46").Value = 207
$ws.Range("F49").Value = 14
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 34
$ws.Range("F9").Value = 14
$ws.Range("F13").Value = 159
$ws.Range("F23").Value = 929
$ws.Range("F24").Value = 14
$ws.Range("F25").Value = 1030
$ws.Range("F26").Value = 221
$ws.Range("F27").Value = 637
$ws.Range("F29").Value = 208
$ws.Range("F30").Value = 3
$ws.Range("F32").Value = 15
$ws.Range("F34").Value = 114
$ws.Range("F35").Value = 157
$ws.Range("F39").Value = 96
$ws.Range("F40").Value = 3
$ws.Range("F41").Value = 19
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 873
$ws.Range("F6").Value = 299
$ws.Range("F7").Value = 142
$ws.Range("F8").Value = 2045
$ws.Range("F9").Value = 3072
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1584
$ws.Range("F4").Value = 121
$ws.Range("F6").Value = 9043
$ws.Range("F7").Value = 299
$ws.Range("F8").Value = 142
$ws.Range("F10").Value = 258
$ws.Range("F11").Value = 2045
$ws.Range("F12").Value = 3072
$ws.Range("F14").Value = 1262
$ws.Range("F15").Value = 159
$ws.Range("F16").Value = 586
$ws.Range("F17").Value = 598
$ws.Range("F18").Value = 289
$ws.Range("F19").Value = 53
$ws.Range("F20").Value = 1481
$ws.Range("F21").Value = 1317
$ws.Range("F23").Value = 1364
$ws.Range("F24").Value = 77
$ws.Range("F26").Value = 86
$ws.Range("F27").Value = 51
$ws.Range("F28").Value = 298
$ws.Range("F29").Value = 1065
$ws.Range("F32").Value = 225
$ws.Range("F33").Value = 14
$ws.Range("F34").Value = 221
$ws.Range("F38").Value = 637
$ws.Range("F40").Value = 208
$ws.Range("F42").Value = 487
$ws.Range("F43").Value = 683
$ws.Range("F44").Value = 114
$ws.Range("F46").Value = 207
$ws.Range("F48").Value = 96
$ws.Range("F50").Value = 19
